$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("G2").Value = -0.0400523560209424
$ws.Range("H2").Value = -0.06858638743455497
$ws.Range("I2").Value = -0.1230366492146597
$ws.Range("J2").Value = -0.1230366492146597
$ws.Range("K2").Value = -3.19
$ws.Range("L2").Value = -0.1670157068062827
$ws.Range("U2").Value = 5.11
$ws.Range("V2").Value = 0.08295454545454546
$ws.Range("W2").Value = -0.1945121951219512
$ws.Range("X2").Value = 0.06303818991073366
$ws.Range("Y2").Value = -0.2575503850326849
$ws.Range("Z2").Value = 5.900525177633615
$ws.Range("AA2").Value = -0.7259808464627746
$ws.Range("AB2").Value = 0.06258170218835768
$ws.Range("AC2").Value = -0.7885625486511323
$ws.Range("AD2").Value = 0.639
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0.639
$ws.Range("AG2").Value = -4.471
$ws.Range("AH2").Value = 0.01026687446777744
$ws.Range("AI2").Value = 0.0430622009569378
$ws.Range("AJ2").Value = -0.07826147840851407
$ws.Range("AK2").Value = -0.4595539109877686
$ws.Range("AL2").Value = 1.8
$ws.Range("AM2").Value = 1.8
$ws.Range("AN2").Value = -0.4915384615384615
$ws.Range("AO2").Value = -1.305555555555556
$ws.Range("AP2").Value = 3.439230769230769
$ws.Range("AQ2").Value = -1.305555555555556
$ws.Range("G3").Value = -0.0400523560209424
$ws.Range("H3").Value = -0.06858638743455497
$ws.Range("I3").Value = -0.1230366492146597
$ws.Range("J3").Value = -0.1230366492146597
$ws.Range("K3").Value = -3.19
$ws.Range("L3").Value = -0.1670157068062827
$ws.Range("U3").Value = 5.11
$ws.Range("V3").Value = 0.08295454545454546
$ws.Range("W3").Value = -0.1945121951219512
$ws.Range("X3").Value = 0.06303818991073366
$ws.Range("Y3").Value = -0.2575503850326849
$ws.Range("Z3").Value = 5.900525177633615
$ws.Range("AA3").Value = -0.7259808464627746
$ws.Range("AB3").Value = 0.06258170218835768
$ws.Range("AC3").Value = -0.7885625486511323
$ws.Range("AD3").Value = 0.639
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.639
$ws.Range("AG3").Value = -4.471
$ws.Range("AH3").Value = 0.01026687446777744
$ws.Range("AI3").Value = 0.0430622009569378
$ws.Range("AJ3").Value = -0.07826147840851407
$ws.Range("AK3").Value = -0.4595539109877686
$ws.Range("AL3").Value = 1.8
$ws.Range("AM3").Value = 1.8
$ws.Range("AN3").Value = -0.4915384615384615
$ws.Range("AO3").Value = -1.305555555555556
$ws.Range("AP3").Value = 3.439230769230769
$ws.Range("AQ3").Value = -1.305555555555556
